{"js": "const replacements = [\n  [\"2026-01-04 Sunday\", \"2026-01-05 Monday\"],\n  [\"822\u00f74=\", \"310\u00f78=\"],\n  [\"626\u00f78=\", \"310\u00f72=\"],\n  [\"795\u00f75=\", \"914\u00f72=\"],\n  [\"284\u00f73=\", \"820\u00f72=\"],\n  [\"997\u00f77=\", \"642\u00f79=\"],\n  [\"211\u00f73=\", \"514\u00f79=\"],\n  [\"971\u00f75=\", \"570\u00f72=\"],\n  [\"854\u00f78=\", \"999\u00f78=\"],\n  [\"910\u00f75=\", \"584\u00f73=\"],\n  [\"270\u00f78=\", \"577\u00f72=\"],\n  [\"207\u00f75=\", \"107\u00f73=\"],\n  [\"762\u00f73=\", \"397\u00f77=\"],\n  [\"166\u00f72=\", \"606\u00f77=\"],\n  [\"894\u00f74=\", \"247\u00f79=\"],\n  [\"696\u00f72=\", \"484\u00f77=\"],\n  [\"247\u00f75=\", \"603\u00f72=\"],\n  [\"110\u00f76=\", \"648\u00f75=\"],\n  [\"281\u00f74=\", \"642\u00f77=\"],\n  [\"595\u00f74=\", \"166\u00f79=\"],\n  [\"129\u00f74=\", \"345\u00f79=\"],\n  [\"301\u00f73=\", \"553\u00f77=\"],\n  [\"716\u00f72=\", \"981\u00f73=\"],\n  [\"426\u00f79=\", \"445\u00f73=\"],\n  [\"182\u00f74=\", \"175\u00f77=\"],\n  [\"296\u00f78=\", \"388\u00f73=\"],\n];\n\nconst docBody = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = docBody.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = '2026-01-04 Sunday'; New = '2026-01-05 Monday'},\n    @{Old = '822\u00f74='; New = '310\u00f78='},\n    @{Old = '626\u00f78='; New = '310\u00f72='},\n    @{Old = '795\u00f75='; New = '914\u00f72='},\n    @{Old = '284\u00f73='; New = '820\u00f72='},\n    @{Old = '997\u00f77='; New = '642\u00f79='},\n    @{Old = '211\u00f73='; New = '514\u00f79='},\n    @{Old = '971\u00f75='; New = '570\u00f72='},\n    @{Old = '854\u00f78='; New = '999\u00f78='},\n    @{Old = '910\u00f75='; New = '584\u00f73='},\n    @{Old = '270\u00f78='; New = '577\u00f72='},\n    @{Old = '207\u00f75='; New = '107\u00f73='},\n    @{Old = '762\u00f73='; New = '397\u00f77='},\n    @{Old = '166\u00f72='; New = '606\u00f77='},\n    @{Old = '894\u00f74='; New = '247\u00f79='},\n    @{Old = '696\u00f72='; New = '484\u00f77='},\n    @{Old = '247\u00f75='; New = '603\u00f72='},\n    @{Old = '110\u00f76='; New = '648\u00f75='},\n    @{Old = '281\u00f74='; New = '642\u00f77='},\n    @{Old = '595\u00f74='; New = '166\u00f79='},\n    @{Old = '129\u00f74='; New = '345\u00f79='},\n    @{Old = '301\u00f73='; New = '553\u00f77='},\n    @{Old = '716\u00f72='; New = '981\u00f73='},\n    @{Old = '426\u00f79='; New = '445\u00f73='},\n    @{Old = '182\u00f74='; New = '175\u00f77='},\n    @{Old = '296\u00f78='; New = '388\u00f73='},\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: text not found -> $($r.Old)\"\n    }\n}\n"}
